$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two more sequential values (14, 15) in P1:Q1,
# using the same style/formatting as the rest of the header row (O1).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# For each data row (2-25), update columns I, K, M, O and add new columns P, Q.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new) = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q (new) = 2
}
